# Append a new row 3 to the "David Miller " sheet, duplicating row 2's
# per-match stats (same venue/date/result/teams/batsman, all counting
# stats at 0) as a second scraped row - matches the source diff which
# grows the sheet from A1:K2 to A1:K3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow = 2
$targetRow = 3

# Columns G:J ("totalRuns","totalBalls","total4s","total6s") hold digit
# strings (e.g. "0") that must stay text, like the rest of the sheet -
# format the destination cells as text *before* assigning so Excel
# doesn't auto-coerce them to numbers.
$numericLookingCols = @("G", "H", "I", "J")

foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")) {
    $src = $ws.Range($col + $sourceRow)
    $dst = $ws.Range($col + $targetRow)

    if ($numericLookingCols -contains $col) {
        $dst.NumberFormat = "@"
    }

    $dst.Value = $src.Value2
}
